$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 86, shifting existing rows 86:150 down to 87:151
$ws.Rows.Item(86).Insert()

$newRow = 86
$ws.Cells.Item($newRow, 1).Value = 5
$ws.Cells.Item($newRow, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item($newRow, 3).Value = "Maule"
$ws.Cells.Item($newRow, 4).Value = 44830
$ws.Cells.Item($newRow, 5).Value = 7
$ws.Cells.Item($newRow, 6).Value = "Fruta"
$ws.Cells.Item($newRow, 7).Value = 100108
$ws.Cells.Item($newRow, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item($newRow, 9).Value = 100108002
$ws.Cells.Item($newRow, 10).Value = "Mango"
$ws.Cells.Item($newRow, 11).Value = "Sin especificar"
$ws.Cells.Item($newRow, 12).Value = "Primera"
$ws.Cells.Item($newRow, 13).Value = 228
$ws.Cells.Item($newRow, 14).Value = 8000
$ws.Cells.Item($newRow, 15).Value = 8000
$ws.Cells.Item($newRow, 16).Value = 8000
$ws.Cells.Item($newRow, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item($newRow, 18).Value = "Brasil"
$ws.Cells.Item($newRow, 19).Value = 2000
$ws.Cells.Item($newRow, 20).Value = 4
